# update MIN / MAX values to align with MHV-17222
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

# Version: 0.2.9-beta -> 0.2.10-beta
$ws.Range("B3").Value = "0.2.10-beta"

# Date: 2023-02-16T09:21:54-06:00 -> 2023-12-06T12:46:33-06:00
$ws.Range("B8").Value = "2023-12-06T12:46:33-06:00"
